# "la til force recomput i load_or_update"
#
# The upstream loader re-ran with the new "force recompute" behaviour, which:
#  - dropped the two earlier placeholder/duplicate "nowind" rows (the ones
#    without amp/freq/per metadata) that used to sit at rows 2 and 3, letting
#    every later row shift up by two, and
#  - recomputed the `ak` column for the row that ends up at row 3 (the
#    "fullpanel-nowind ... run2" row), which came out with a tiny
#    floating-point-level difference from before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old row 2 ("fullpanel-nowind-NO-...run2") and, after the shift,
# the new row 2 ("fullpanel-nowind-...run2" without amp/freq/per) — this
# moves every following row up by two, matching the new A1:AZ7 dimension.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Force-recompute artifact: the "ak" value for the
# fullpanel-nowind-amp0100-freq1300-per30-...-run2 row (now row 3) comes back
# very slightly different after the recompute.
$ws.Range("AO3").Value = 0.05873700892101653
